$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A2").Value = 111396045
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = '70'
$ws.Range("J2").Value = 'stjälkar/strån/skott'
$ws.Range("P2").Value = 'S om järnvägen, Vg'
$ws.Range("Q2").Value = 431889.3909100805
$ws.Range("R2").Value = 6419670.266848063
$ws.Range("A3").Value = 111396060
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = '90'
$ws.Range("J3").ClearContents()
$ws.Range("P3").Value = 'S om järnvägen - 3, Vg'
$ws.Range("Q3").Value = 432076.641898193
$ws.Range("R3").Value = 6419661.774153749
$ws.Range("A5").Value = 111482955
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = '70'
$ws.Range("P5").Value = 'S om järnvägen - 5, Vg'
$ws.Range("Q5").Value = 432064.1298546481
$ws.Range("R5").Value = 6419677.395781181
$ws.Range("A6").Value = 111482936
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = '25'
$ws.Range("P6").Value = 'S om järnvägen - 4, Vg'
$ws.Range("Q6").Value = 432073.5656663703
$ws.Range("R6").Value = 6419668.734013095
$ws.Range("A7").Value = 111482980
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = '10'
$ws.Range("P7").Value = 'S om järnvägen - 6, Vg'
$ws.Range("Q7").Value = 432048.2263952638
$ws.Range("R7").Value = 6419681.385014677
$ws.Range("A8").Value = 111483107
$ws.Range("B8").Value = 73681
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 6439
$ws.Range("F8").Value = 'Gulnål'
$ws.Range("G8").Value = 'Chaenotheca brachypoda'
$ws.Range("H8").Value = '(Ach.) Tibell'
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("P8").Value = 'S om järnvägen - 8, Vg'
$ws.Range("Q8").Value = 431947.1499479365
$ws.Range("R8").Value = 6419623.056550305
$ws.Range("AJ8").Value = 'tall'
$ws.Range("AK8").Value = 'Pinus sylvestris'
$ws.Range("AM8").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO8").Value = 'Standing dead tree/snags # Pinus sylvestris'
$ws.Range("A9").Value = 111483462
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = '45'
$ws.Range("K9").Value = 'blomning'
$ws.Range("P9").Value = 'S om järnvägen - 16, Vg'
$ws.Range("Q9").Value = 431654.0242198514
$ws.Range("R9").Value = 6419791.70470859
$ws.Range("A10").Value = 111483140
$ws.Range("B10").Value = 73683
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 306
$ws.Range("F10").Value = 'Kornig nållav'
$ws.Range("G10").Value = 'Chaenotheca chlorella'
$ws.Range("H10").Value = '(Ach.) Müll.Arg.'
$ws.Range("P10").Value = 'S om järnvägen - 9, Vg'
$ws.Range("Q10").Value = 431942.9372677525
$ws.Range("R10").Value = 6419625.784949708
$ws.Range("A11").Value = 111483437
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = '100'
$ws.Range("P11").Value = 'S om järnvägen - 15, Vg'
$ws.Range("Q11").Value = 431797.479853621
$ws.Range("R11").Value = 6419681.394993878
$ws.Range("A12").Value = 111483105
$ws.Range("B12").Value = 73689
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 308
$ws.Range("F12").Value = 'Brunpudrad nållav'
$ws.Range("G12").Value = 'Chaenotheca gracillima'
$ws.Range("H12").Value = '(Vain.) Tibell'
$ws.Range("I12").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("P12").Value = 'S om järnvägen - 8, Vg'
$ws.Range("Q12").Value = 431947.1499479365
$ws.Range("R12").Value = 6419623.056550305
$ws.Range("AJ12").Value = 'tall'
$ws.Range("AK12").Value = 'Pinus sylvestris'
$ws.Range("AM12").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO12").Value = 'Standing dead tree/snags # Pinus sylvestris'
$ws.Range("A13").Value = 111483037
$ws.Range("B13").Value = 96348
$ws.Range("D13").Value = 'VU'
$ws.Range("E13").Value = 220787
$ws.Range("F13").Value = 'Knärot'
$ws.Range("G13").Value = 'Goodyera repens'
$ws.Range("H13").Value = '(L.) R. Br.'
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = '60'
$ws.Range("J13").Value = 'stjälkar/strån/skott'
$ws.Range("K13").Value = 'blomning'
$ws.Range("L13").ClearContents()
$ws.Range("P13").Value = 'S om järnvägen - 7, Vg'
$ws.Range("Q13").Value = 432060.6482816387
$ws.Range("R13").Value = 6419660.45125766
$ws.Range("AJ13").ClearContents()
$ws.Range("AK13").ClearContents()
$ws.Range("AM13").ClearContents()
$ws.Range("AO13").ClearContents()
$ws.Range("A14").Value = 111490843
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = '50'
$ws.Range("K14").Value = 'fullt utvecklade blad'
$ws.Range("P14").Value = 'S om järnvägen - 17, Vg'
$ws.Range("Q14").Value = 431803.2980747336
$ws.Range("R14").Value = 6419679.170503675
$ws.Range("A15").Value = 111491187
$ws.Range("B15").Value = 96348
$ws.Range("D15").Value = 'VU'
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = 'Knärot'
$ws.Range("G15").Value = 'Goodyera repens'
$ws.Range("H15").Value = '(L.) R. Br.'
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = '60'
$ws.Range("J15").Value = 'stjälkar/strån/skott'
$ws.Range("K15").Value = 'blomning'
$ws.Range("L15").ClearContents()
$ws.Range("P15").Value = 'S om järnvägen - 18, Vg'
$ws.Range("Q15").Value = 431829.514510141
$ws.Range("R15").Value = 6419749.394753682
$ws.Range("AJ15").ClearContents()
$ws.Range("AK15").ClearContents()
$ws.Range("AM15").ClearContents()
$ws.Range("AO15").ClearContents()
$ws.Range("A16").Value = 111491635
$ws.Range("B16").Value = 96348
$ws.Range("D16").Value = 'VU'
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = 'Knärot'
$ws.Range("G16").Value = 'Goodyera repens'
$ws.Range("H16").Value = '(L.) R. Br.'
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = '10'
$ws.Range("J16").Value = 'stjälkar/strån/skott'
$ws.Range("K16").Value = 'blomning'
$ws.Range("L16").ClearContents()
$ws.Range("P16").Value = 'S om järnvägen - 21, Vg'
$ws.Range("Q16").Value = 431859.6228004749
$ws.Range("R16").Value = 6419672.898494411
$ws.Range("AM16").ClearContents()
$ws.Range("AO16").ClearContents()
$ws.Range("A17").Value = 111483300
$ws.Range("B17").Value = 73689
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 308
$ws.Range("F17").Value = 'Brunpudrad nållav'
$ws.Range("G17").Value = 'Chaenotheca gracillima'
$ws.Range("H17").Value = '(Vain.) Tibell'
$ws.Range("I17").ClearContents()
$ws.Range("J17").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("P17").Value = 'S om järnvägen - 12, Vg'
$ws.Range("Q17").Value = 431888.091041417
$ws.Range("R17").Value = 6419625.122914318
$ws.Range("AJ17").Value = 'tall'
$ws.Range("AK17").Value = 'Pinus sylvestris'
$ws.Range("AM17").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO17").Value = 'Standing dead tree/snags # Pinus sylvestris'
$ws.Range("A18").Value = 111483197
$ws.Range("P18").Value = 'S om järnvägen - 11, Vg'
$ws.Range("Q18").Value = 431937.082796899
$ws.Range("R18").Value = 6419625.884406033
$ws.Range("AJ18").ClearContents()
$ws.Range("AK18").ClearContents()
$ws.Range("AO18").Value = 'Standing dead tree/snags'
$ws.Range("A19").Value = 111483381
$ws.Range("B19").Value = 73689
$ws.Range("D19").Value = 'NT'
$ws.Range("E19").Value = 308
$ws.Range("F19").Value = 'Brunpudrad nållav'
$ws.Range("G19").Value = 'Chaenotheca gracillima'
$ws.Range("H19").Value = '(Vain.) Tibell'
$ws.Range("P19").Value = 'S om järnvägen - 14, Vg'
$ws.Range("Q19").Value = 431754.10213514
$ws.Range("R19").Value = 6419728.893211351
$ws.Range("A20").Value = 111661832
$ws.Range("Q20").Value = 432076.4609239195
$ws.Range("R20").Value = 6419682.500295377
$ws.Range("A22").Value = 111661838
$ws.Range("Q22").Value = 431799.2483237319
$ws.Range("R22").Value = 6419691.460736625
